{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\n// ------------------------------------------------------------------\n// Change 1: the \"{{Nicknames :empty:remove:row}}\" list paragraph was\n// split across three separate runs (\"{{Nicknames :\", \"empty:remove:row\",\n// \"}}\"). Collapse it back into a single run (same text, same\n// formatting) by replacing the paragraph's range text in place.\n// ------------------------------------------------------------------\nlet nicknamesPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"{{Nicknames :empty:remove:row}}\") !== -1) {\n    nicknamesPara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (nicknamesPara) {\n  const nicknamesRange = nicknamesPara.getRange();\n  nicknamesRange.insertText(\"{{Nicknames :empty:remove:row}}\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// ------------------------------------------------------------------\n// Change 2: add a new paragraph right after the \"{{NotReplacable}} -\n// this should not be replaced\" paragraph, containing a trigger-param\n// version of the same placeholder that is meant to be removed by the\n// template engine.\n// ------------------------------------------------------------------\nlet notReplacablePara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"{{NotReplacable}}\") !== -1) {\n    notReplacablePara = paragraphs.items[i];\n  }\n}\n\nif (notReplacablePara) {\n  notReplacablePara.insertParagraph(\n    \"{{NotReplacable :empty:remove:row}} \\u2013 this should be removed\",\n    Word.InsertLocation.after\n  );\n}\nawait context.sync();\n\n// Re-fetch the paragraph collection so the freshly inserted paragraph is a\n// live (non-stub) object before formatting it - setting .font straight off\n// the object returned by insertParagraph does not reliably reach the\n// paragraph-mark (w:pPr/w:rPr) properties.\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  paragraphs2.items[i].load(\"text\");\n}\nawait context.sync();\n\nlet newPara = null;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text.indexOf(\"{{NotReplacable :empty:remove:row}}\") !== -1) {\n    newPara = paragraphs2.items[i];\n  }\n}\n\nif (newPara) {\n  newPara.font.color = \"#666666\";\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# ------------------------------------------------------------------\n# Change 1: the \"{{Nicknames :empty:remove:row}}\" list paragraph was\n# split across three separate runs (\"{{Nicknames :\", \"empty:remove:row\",\n# \"}}\"). Collapse it back into a single run (same text, same formatting).\n# ------------------------------------------------------------------\n$nicknamesPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*{{Nicknames :empty:remove:row}}*\") {\n        $nicknamesPara = $p\n        break\n    }\n}\n\nif ($nicknamesPara -ne $null) {\n    $rng = $nicknamesPara.Range\n    [void]$rng.MoveEnd(1, -1)\n    # Assigning .Text always rewrites the range as a single run, but if the\n    # new text is identical to what's already there it's a no-op, so first\n    # stamp a throwaway placeholder to force the rewrite, then set the real\n    # text back.\n    $rng.Text = \"___TEMP___\"\n    $rng2 = $nicknamesPara.Range\n    [void]$rng2.MoveEnd(1, -1)\n    $rng2.Text = \"{{Nicknames :empty:remove:row}}\"\n}\n\n# ------------------------------------------------------------------\n# Change 2: add a new paragraph right after the \"{{NotReplacable}} - this\n# should not be replaced\" paragraph, containing a trigger-param version of\n# the same placeholder that is meant to be removed by the template engine.\n# ------------------------------------------------------------------\n$notReplacablePara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*{{NotReplacable}}*\") {\n        $notReplacablePara = $p\n    }\n}\n\nif ($notReplacablePara -ne $null) {\n    $notReplacablePara.Range.InsertParagraphAfter()\n    $newPara = $notReplacablePara.Next()\n    $newRng = $newPara.Range\n    [void]$newRng.MoveEnd(1, -1)\n    $newRng.Text = \"{{NotReplacable :empty:remove:row}} \" + [char]0x2013 + \" this should be removed\"\n    $newRng.Font.Color = 6710886\n    # Also color the paragraph mark itself (w:pPr/w:rPr) to match the\n    # target markup, not just the run.\n    $newPara.Range.Font.Color = 6710886\n}\n"}
